$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 15 cell values ---
$ws.Range("B15").Value = "taxi game"
$ws.Range("C15").Value = "irisalmog47@gmail.com"
$ws.Range("D15").Value = "bittonnir12@gmail.com"
$ws.Range("F15").Value = "I made it to the final level. I cant believe it. Greatest taxi game"

# --- Row height for row 15 (matches Excel's autofit after edit) ---
$ws.Rows.Item(15).RowHeight = 13.8

# --- Rebuild hyperlinks: the runtime's Hyperlinks.Delete() on any range
#     clears the whole sheet's collection, so remove them all then
#     re-add every link except the one we dropped (old C15, which had
#     no hyperlink in the new content), with D15 pointing at its new
#     address. ---
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:snizzvered@gmail.com", "", "", "snizzvered@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:krigelron@gmail.com", "", "", "krigelron@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:jorjkluni03@gmail.com", "", "", "jorjkluni03@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:sugarderryfireapp@gmail.com", "", "", "sugarderryfireapp@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:sugarderryfire@gmail.com", "", "", "sugarderryfire@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:sm6502345@gmail.com", "", "", "sm6502345@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:cybworking@gmail.com", "", "", "cybworking@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:eligitel@gmail.com", "", "", "eligitel@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:ronenchen27@gmail.com", "", "", "ronenchen27@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:sixsevensix67676@gmail.com", "", "", "sixsevensix67676@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:dony1098765432@gmail.com", "", "", "dony1098765432@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:sixsevensix67676@gmail.com", "", "", "sixsevensix67676@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:gregneri12@gmail.com", "", "", "gregneri12@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:halachme@gmail.com", "", "", "halachme@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:nitanoren23@gmail.com", "", "", "nitanoren23@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:cristianjohn1222@gmail.com", "", "", "cristianjohn1222@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:bittonnir12@gmail.com", "", "", "bittonnir12@gmail.com") | Out-Null
